$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry row (row 33) - "Finalizacion" / "Fin de proyecto"
$ws.Range("C33").Value = "Finalizacion"
$ws.Range("D33").Value = (Get-Date -Year 2022 -Month 7 -Day 4 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("D33").NumberFormat = "d-mmm"
$ws.Range("E33").Value = "Fin de proyecto"

# Match the scrolled view / selection left behind by the edit
[void]$ws.Range("E33:E34").Select()
$excel.ActiveWindow.ScrollRow = 25
